# semana 37 de 2024
# Updates the Esperado/Observado/valor p columns (C/D/E) for the existing
# events and inserts three new event rows (610, 720, 730) that appear in
# this week's report, shifting the remaining rows down.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Refresh C/D/E for the events that keep their current row ----------
# row -> (Esperado, Observado, "valor p")
$updates = @{
    3  = @(1, 5, 0)
    4  = @(0, 0, 1)
    5  = @(6, 15, 0)
    6  = @(1, 69, 0)
    7  = @(2, 3, 0.18)
    9  = @(0, 1, 0)
    11 = @(43, 33, 0.02)
    12 = @(0, 0, 1)
    13 = @(1, 1, 0.37)
    14 = @(4, 10, 0.01)
    15 = @(1, 0, 0.37)
    16 = @(0, 0, 1)
    17 = @(1, 0, 0.37)
    18 = @(0, 0, 1)
    19 = @(10, 8, 0.11)
    21 = @(6, 4, 0.13)
    24 = @(0, 0, 1)
    26 = @(0, 1, 0)
    27 = @(8, 2, 0.01)
    28 = @(1, 1, 0.37)
    29 = @(0, 0, 1)
}

foreach ($r in $updates.Keys) {
    $vals = $updates[$r]
    $ws.Cells.Item($r, 3).Value = $vals[0]
    $ws.Cells.Item($r, 4).Value = $vals[1]
    $ws.Cells.Item($r, 5).Value = $vals[2]
}

# --- 2. Insert three new event rows (610, 720, 730) -------------------------
# They are not contiguous with the old data: "610" lands right before
# "620 Parotiditis" (old row 30), and "720"/"730" land right after it, before
# "740 Sifilis congenita" (old row 31). Insert bottom-up so row numbers used
# below stay valid as each insert shifts everything after it down.

# 2a. Two rows between old row 30 (Parotiditis) and old row 31 (Sifilis congenita)
$ws.Range("A31:A32").EntireRow.Insert()

# 2b. One row before old row 30 (Parotiditis)
$ws.Range("A30:A30").EntireRow.Insert()

# --- 3. Write the new events into the freshly inserted rows -----------------
# The "evento" column holds codes that look numeric ("610", "720", ...) but
# are stored as text throughout the sheet, so force a text format before
# assigning the value (otherwise Excel auto-detects them as numbers).
$ws.Cells.Item(30, 1).NumberFormat = "@"
$ws.Cells.Item(30, 1).Value = "610"
$ws.Cells.Item(30, 4).Value = 0

$ws.Cells.Item(32, 1).NumberFormat = "@"
$ws.Cells.Item(32, 1).Value = "720"
$ws.Cells.Item(32, 2).Value = "Sindrome de rubeola congenita"
$ws.Cells.Item(32, 3).Value = 0
$ws.Cells.Item(32, 4).Value = 0
$ws.Cells.Item(32, 5).Value = 1

$ws.Cells.Item(33, 1).NumberFormat = "@"
$ws.Cells.Item(33, 1).Value = "730"
$ws.Cells.Item(33, 2).Value = "Sarampion"
$ws.Cells.Item(33, 3).Value = 0
$ws.Cells.Item(33, 4).Value = 1
$ws.Cells.Item(33, 5).Value = 0

# --- 4. Refresh C/D/E for the events that shifted down with the inserts ----
# (row numbers below are the NEW row numbers after both inserts)
$shiftedUpdates = @{
    31 = @(0, 0, 1)     # 620 Parotiditis (was row 30)
    34 = @(0, 1, 0)     # 740 Sifilis congenita (was row 31)
    35 = @(3, 2, 0.22)  # 750 Sifilis gestacional (was row 32)
    36 = @(7, 4, 0.09)  # 813 Tuberculosis (was row 33)
    37 = @(11, 2, 0)    # 831 Varicela individual (was row 34)
}

foreach ($r in $shiftedUpdates.Keys) {
    $vals = $shiftedUpdates[$r]
    $ws.Cells.Item($r, 3).Value = $vals[0]
    $ws.Cells.Item($r, 4).Value = $vals[1]
    $ws.Cells.Item($r, 5).Value = $vals[2]
}

# Rows 38 (850 Vih/sida, was row 35) and 39 (895 Zika, was row 36) keep the
# same C/D/E values they had before the inserts, so nothing else to change.
